# "completed part 3 and 4" — the old "Audio"/audio-file column (column A,
# containing the header "Audio" and a merged A2:A4 cell with "audio.mp3")
# is removed entirely. Every other column shifts one position to the left:
#   B..G (answer choices, "Dap an" letter, "Cau hoi" question) -> A..F
#
# Deleting the whole column (rather than just clearing+rewriting cells)
# is what reproduces the shift of dimension/cols/row spans/selection seen
# in the target worksheet, and lets the engine drop the now-unreferenced
# "Audio"/"audio.mp3" shared strings automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns.Item(1).Delete()

# Matches the post-edit selection state recorded in the sheet view
# (the just-deleted column's position, now occupied by the shifted data,
# ends up fully selected as column A).
$ws.Range("A1:A1048576").Select()
